# Auto-generated edit script applying cell text updates from the target diff.
# Cells D/E (and B/C for the swapped rows) hold plain text values (numbers-as-text,
# percentages, names, links). We force the cell to Text format before assigning so
# Excel does not reinterpret strings like "224.00" or "1.705.71" as numeric/date values,
# then restore the cell style to "Normal" so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '27.324.70'
Set-TextValue 'E2' '  -0.78%  '
Set-TextValue 'D3' '1.706.45'
Set-TextValue 'E3' '  -1.06%  '
Set-TextValue 'E4' '  -0.09%  '
Set-TextValue 'D5' '224.00'
Set-TextValue 'E5' '  -0.83%  '
Set-TextValue 'D6' '0.5323'
Set-TextValue 'E6' '  -1.07%  '
Set-TextValue 'E7' '  -0.06%  '
Set-TextValue 'D8' '0.2667'
Set-TextValue 'E8' '  -0.27%  '
Set-TextValue 'D9' '0.06610'
Set-TextValue 'E9' '  +0.00%  '
Set-TextValue 'E10' '  -4.63%  '
Set-TextValue 'D11' '0.07652'
Set-TextValue 'E11' '  -1.07%  '
Set-TextValue 'E12' '  -2.31%  '
Set-TextValue 'D13' '1.941.20'
Set-TextValue 'E13' '  -1.04%  '
Set-TextValue 'D14' '1.706.08'
Set-TextValue 'E14' '  -1.06%  '
Set-TextValue 'D15' '0.5824'
Set-TextValue 'E15' '  -0.96%  '
Set-TextValue 'E16' '  -1.72%  '
Set-TextValue 'D17' '67.69'
Set-TextValue 'E17' '  -0.55%  '
Set-TextValue 'D18' '27.332.19'
Set-TextValue 'D19' '216.05'
Set-TextValue 'E19' '  -2.71%  '
Set-TextValue 'E20' '  -0.05%  '
Set-TextValue 'D21' '4.638'
Set-TextValue 'E21' '  -2.33%  '
Set-TextValue 'E22' '  -2.72%  '
Set-TextValue 'D23' '5.990'
Set-TextValue 'E23' '  -1.84%  '
Set-TextValue 'E24' '  -0.09%  '
Set-TextValue 'D25' '143.65'
Set-TextValue 'E25' '  -3.12%  '
Set-TextValue 'D26' '1.703'
Set-TextValue 'E26' '  +0.44%  '
Set-TextValue 'D27' '0.1204'
Set-TextValue 'E27' '  -2.37%  '
Set-TextValue 'D28' '7.227'
Set-TextValue 'E28' '  -2.37%  '
Set-TextValue 'D29' '16.21'
Set-TextValue 'E29' '  -2.82%  '
Set-TextValue 'D30' '0.05373'
Set-TextValue 'E30' '  -2.93%  '
Set-TextValue 'D32' '3.480'
Set-TextValue 'E32' '  -1.86%  '
Set-TextValue 'D33' '3.418'
Set-TextValue 'E33' '  -1.31%  '
Set-TextValue 'D34' '1.647'
Set-TextValue 'E34' '  -1.00%  '
Set-TextValue 'D35' '2.863'
Set-TextValue 'E35' '  +1.46%  '
Set-TextValue 'D36' '0.9506'
Set-TextValue 'E36' '  -1.20%  '
Set-TextValue 'D37' '2.402'
Set-TextValue 'E37' '  -1.80%  '
Set-TextValue 'D38' '0.5854'
Set-TextValue 'E38' '  -1.74%  '
Set-TextValue 'D39' '0.01640'
Set-TextValue 'E39' '  -0.41%  '
Set-TextValue 'D40' '5.811'
Set-TextValue 'E40' '  -1.98%  '
Set-TextValue 'D41' '1.044.30'
Set-TextValue 'E41' '  -1.50%  '
Set-TextValue 'B42' 'TrustWalletToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D42' '0.8433'
Set-TextValue 'E42' '  -1.23%  '
Set-TextValue 'B43' 'PaxDollar'
Set-TextValue 'C43' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D43' '1.003'
Set-TextValue 'E43' '  -0.01%  '
Set-TextValue 'D44' '100.90'
Set-TextValue 'E44' '  -0.71%  '
Set-TextValue 'D45' '1.848.79'
Set-TextValue 'E45' '  -1.03%  '
Set-TextValue 'E46' '  -4.38%  '
Set-TextValue 'D47' '57.89'
Set-TextValue 'E47' '  -2.07%  '
Set-TextValue 'D48' '0.4523'
Set-TextValue 'E48' '  +1.89%  '
Set-TextValue 'E49' '  +0.18%  '
Set-TextValue 'D50' '8.080'
Set-TextValue 'E50' '  -1.76%  '
Set-TextValue 'D51' '0.05231'
Set-TextValue 'E51' '  -0.90%  '
